# Auto-generated edit script applying numeric corrections to H:N columns
# across multiple sheets, per the supplied OOXML diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 669.6957
$ws.Range("I28").Value = 386.2
$ws.Range("K28").Value = 386.2
$ws.Range("M28").Value = 98.80000000000001
$ws.Range("H62").Value = 3262.5
$ws.Range("I62").Value = 2915
$ws.Range("K62").Value = 2915
$ws.Range("M62").Value = -2291
$ws.Range("H65").Value = 3262.5
$ws.Range("I65").Value = 2915
$ws.Range("K65").Value = 14575
$ws.Range("M65").Value = -11455
$ws.Range("H74").Value = 3320.8572
$ws.Range("I74").Value = 3200.75
$ws.Range("J74").Value = 3368.9
$ws.Range("K74").Value = 3200.75
$ws.Range("L74").Value = 3368.9
$ws.Range("M74").Value = -2264.75
$ws.Range("N74").Value = -5240.9
$ws.Range("H77").Value = 3320.8572
$ws.Range("I77").Value = 3200.75
$ws.Range("J77").Value = 3368.9
$ws.Range("K77").Value = 16003.75
$ws.Range("L77").Value = 16844.5
$ws.Range("M77").Value = -11323.75
$ws.Range("N77").Value = -26204.5
$ws.Range("H103").Value = 2011.6666
$ws.Range("I103").Value = 2294
$ws.Range("J103").Value = 600
$ws.Range("K103").Value = 6882
$ws.Range("L103").Value = 1800
$ws.Range("M103").Value = -6296
$ws.Range("N103").Value = -2972
$ws.Range("H129").Value = 2806.36
$ws.Range("J129").Value = 901.91895
$ws.Range("L129").Value = 2705.75685
$ws.Range("N129").Value = -12705.75685
$ws.Range("H132").Value = 6256213.5
$ws.Range("I132").Value = 7359454
$ws.Range("J132").Value = 4515.1665
$ws.Range("K132").Value = 22078362
$ws.Range("L132").Value = 13545.4995
$ws.Range("M132").Value = -22075832
$ws.Range("N132").Value = -18605.4995
$ws.Range("H137").Value = 1340.4
$ws.Range("I137").Value = 1306.75
$ws.Range("J137").Value = 1475
$ws.Range("K137").Value = 3920.25
$ws.Range("L137").Value = 4425
$ws.Range("M137").Value = -1370.25
$ws.Range("N137").Value = -9525

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3395.48
$ws.Range("I32").Value = 3109.875
$ws.Range("J32").Value = 10250
$ws.Range("K32").Value = 3109.875
$ws.Range("L32").Value = 10250
$ws.Range("M32").Value = -2822.875
$ws.Range("N32").Value = -10824
$ws.Range("H61").Value = 1359.3135
$ws.Range("I61").Value = 956.6667
$ws.Range("J61").Value = 3031.8462
$ws.Range("K61").Value = 956.6667
$ws.Range("L61").Value = 3031.8462
$ws.Range("M61").Value = -744.6667
$ws.Range("N61").Value = -3455.8462
$ws.Range("H74").Value = 923.18604
$ws.Range("I74").Value = 914.3611
$ws.Range("K74").Value = 914.3611
$ws.Range("M74").Value = -40.36109999999996
$ws.Range("H77").Value = 923.18604
$ws.Range("I77").Value = 914.3611
$ws.Range("K77").Value = 4571.805499999999
$ws.Range("M77").Value = -203.8054999999995
$ws.Range("H97").Value = 34755.5
$ws.Range("I97").Value = 67706
$ws.Range("J97").Value = 1805
$ws.Range("K97").Value = 67706
$ws.Range("L97").Value = 1805
$ws.Range("M97").Value = -67210
$ws.Range("N97").Value = -2797
$ws.Range("H125").Value = 39483.6
$ws.Range("J125").Value = 39483.6
$ws.Range("L125").Value = 39483.6
$ws.Range("N125").Value = -49323.6
$ws.Range("H132").Value = 14237.911
$ws.Range("I132").Value = 15631.743
$ws.Range("J132").Value = 5178
$ws.Range("K132").Value = 46895.229
$ws.Range("L132").Value = 15534
$ws.Range("M132").Value = -44365.229
$ws.Range("N132").Value = -20594
$ws.Range("H136").Value = 1359.3135
$ws.Range("I136").Value = 956.6667
$ws.Range("J136").Value = 3031.8462
$ws.Range("K136").Value = 2870.0001
$ws.Range("L136").Value = 9095.5386
$ws.Range("M136").Value = -320.0001000000002
$ws.Range("N136").Value = -14195.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 11823.667
$ws.Range("I26").Value = 11823.667
$ws.Range("K26").Value = 11823.667
$ws.Range("M26").Value = -11531.667
$ws.Range("H105").Value = 81669.16
$ws.Range("I105").Value = 54313.105
$ws.Range("J105").Value = 168296.67
$ws.Range("K105").Value = 54313.105
$ws.Range("L105").Value = 168296.67
$ws.Range("M105").Value = -52566.105
$ws.Range("N105").Value = -171790.67
$ws.Range("H134").Value = 2371.1292
$ws.Range("I134").Value = 1610.8846
$ws.Range("J134").Value = 6324.4
$ws.Range("K134").Value = 4832.6538
$ws.Range("L134").Value = 18973.2
$ws.Range("M134").Value = -2297.6538
$ws.Range("N134").Value = -24043.2
$ws.Range("H140").Value = 39233
$ws.Range("I140").Value = 20709
$ws.Range("J140").Value = 48495
$ws.Range("K140").Value = 20709
$ws.Range("L140").Value = 48495
$ws.Range("M140").Value = -15529
$ws.Range("N140").Value = -58855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 794.0476
$ws.Range("I107").Value = 782.0714
$ws.Range("J107").Value = 818
$ws.Range("K107").Value = 782.0714
$ws.Range("L107").Value = 818
$ws.Range("M107").Value = 1137.9286
$ws.Range("N107").Value = -4658

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1190.909
$ws.Range("J34").Value = 1290
$ws.Range("L34").Value = 3870
$ws.Range("N34").Value = -4038
$ws.Range("H102").Value = 4821.25
$ws.Range("J102").Value = 4928.3335
$ws.Range("L102").Value = 14785.0005
$ws.Range("N102").Value = -19653.0005
$ws.Range("H120").Value = 448725.4
$ws.Range("I120").Value = 448725.4
$ws.Range("K120").Value = 1346176.2
$ws.Range("M120").Value = -1341338.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 40002010
$ws.Range("I97").Value = 47620884
$ws.Range("K97").Value = 47620884
$ws.Range("M97").Value = -47620388
$ws.Range("H122").Value = 1753.8276
$ws.Range("I122").Value = 1597.0625
$ws.Range("K122").Value = 4791.1875
$ws.Range("M122").Value = -2341.1875
$ws.Range("H132").Value = 2565.738
$ws.Range("I132").Value = 2154
$ws.Range("J132").Value = 3389.2144
$ws.Range("K132").Value = 6462
$ws.Range("L132").Value = 10167.6432
$ws.Range("M132").Value = -3932
$ws.Range("N132").Value = -15227.6432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1490.8695
$ws.Range("I61").Value = 1514.2142
$ws.Range("J61").Value = 1454.5555
$ws.Range("K61").Value = 1514.2142
$ws.Range("L61").Value = 1454.5555
$ws.Range("M61").Value = -1312.2142
$ws.Range("N61").Value = -1858.5555
$ws.Range("H82").Value = 2006.8889
$ws.Range("I82").Value = 1761.4286
$ws.Range("J82").Value = 2163.0908
$ws.Range("K82").Value = 1761.4286
$ws.Range("L82").Value = 2163.0908
$ws.Range("M82").Value = -1400.4286
$ws.Range("N82").Value = -2885.0908
$ws.Range("H85").Value = 2006.8889
$ws.Range("I85").Value = 1761.4286
$ws.Range("J85").Value = 2163.0908
$ws.Range("K85").Value = 1761.4286
$ws.Range("L85").Value = 2163.0908
$ws.Range("M85").Value = -513.4286
$ws.Range("N85").Value = -4659.0908
$ws.Range("H113").Value = 1490.8695
$ws.Range("I113").Value = 1514.2142
$ws.Range("J113").Value = 1454.5555
$ws.Range("K113").Value = 1514.2142
$ws.Range("L113").Value = 1454.5555
$ws.Range("M113").Value = 655.7858000000001
$ws.Range("N113").Value = -5794.5555
$ws.Range("H132").Value = 4621.7915
$ws.Range("I132").Value = 4622
$ws.Range("J132").Value = 4621.4443
$ws.Range("K132").Value = 13866
$ws.Range("L132").Value = 13864.3329
$ws.Range("M132").Value = -11336
$ws.Range("N132").Value = -18924.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1302.8788
$ws.Range("I122").Value = 1302.8788
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3908.6364
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1458.6364
$ws.Range("N122").ClearContents()
